# Append five new measurement rows (Dhandha+25) to the temperature_EoR dataset sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$method = "Cosmic X-ray background, UV luminosity function, 21cm power spectrum"
$reference = "Dhandha+25"

# Row 18
$ws.Range("A18").Value = 15
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = "5-7.7"
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = "6.4-33.9"
$ws.Range("H18").Value = $method
$ws.Range("I18").Value = $reference

# Row 19
$ws.Range("A19").Value = 12.5
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = "3.6-16"
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = "4.5-19"
$ws.Range("H19").Value = $method
$ws.Range("I19").Value = $reference

# Row 20
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = "2.5-66.2"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = "3.1-73.3"
$ws.Range("H20").Value = $method
$ws.Range("I20").Value = $reference

# Row 21
$ws.Range("A21").Value = 8
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = "3.7-349.5"
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = "4.3-359.2"
$ws.Range("H21").Value = $method
$ws.Range("I21").Value = $reference

# Row 22
$ws.Range("A22").Value = 6
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = "19.8-2077.9"
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = "19-1258.1"
$ws.Range("H22").Value = $method
$ws.Range("I22").Value = $reference

# Column D needs to widen slightly to fit the new, longer values.
$ws.Columns.Item(4).ColumnWidth = 10

# Reflect the author's final cursor position/selection when the file was saved.
[void]$ws.Range("D24").Select()
